$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 26.127733
$ws.Range("H2").Value = 78.383199
$ws.Range("I2").Value = 0.2666992864894373
$ws.Range("J2").Value = 0.2666992864894374
$ws.Range("M2").Value = 0.8317113333333332
$ws.Range("N2").Value = 2.495134
$ws.Range("O2").Value = 0.0263454906755698
$ws.Range("P2").Value = 0.0263454906755698
$ws.Range("Q2").Value = 21.73073165040733
$ws.Range("R2").Value = 195.576584853666
$ws.Range("S2").Value = 0.007026323565388588
$ws.Range("T2").Value = 0.007026323565388592
$ws.Range("G3").Value = 26.127733
$ws.Range("H3").Value = 78.383199
$ws.Range("I3").Value = 0.2666992864894373
$ws.Range("J3").Value = 0.2666992864894374
$ws.Range("O3").Value = 0.6529848313028861
$ws.Range("P3").Value = 0.6529848313028862
$ws.Range("Q3").Value = 538.6059540727314
$ws.Range("R3").Value = 4847.453586654582
$ws.Range("S3").Value = 0.1741505885969053
$ws.Range("T3").Value = 0.1741505885969054
$ws.Range("G4").Value = 26.127733
$ws.Range("H4").Value = 78.383199
$ws.Range("I4").Value = 0.2666992864894373
$ws.Range("J4").Value = 0.2666992864894374
$ws.Range("M4").Value = 10.12334933333333
$ws.Range("N4").Value = 30.370048
$ws.Range("O4").Value = 0.3206696780215441
$ws.Range("P4").Value = 0.3206696780215441
$ws.Range("Q4").Value = 264.5001684470614
$ws.Range("R4").Value = 2380.501516023552
$ws.Range("S4").Value = 0.08552237432714339
$ws.Range("T4").Value = 0.08552237432714342
$ws.Range("I5").Value = 0.2440410104700376
$ws.Range("J5").Value = 0.2440410104700377
$ws.Range("M5").Value = 0.8317113333333332
$ws.Range("N5").Value = 2.495134
$ws.Range("O5").Value = 0.0263454906755698
$ws.Range("P5").Value = 0.0263454906755698
$ws.Range("Q5").Value = 19.88452905151911
$ws.Range("R5").Value = 178.960761463672
$ws.Range("S5").Value = 0.006429380165795006
$ws.Range("T5").Value = 0.00642938016579501
$ws.Range("I6").Value = 0.2440410104700376
$ws.Range("J6").Value = 0.2440410104700377
$ws.Range("O6").Value = 0.6529848313028861
$ws.Range("P6").Value = 0.6529848313028862
$ws.Range("S6").Value = 0.1593550780527634
$ws.Range("T6").Value = 0.1593550780527634
$ws.Range("I7").Value = 0.2440410104700376
$ws.Range("J7").Value = 0.2440410104700377
$ws.Range("M7").Value = 10.12334933333333
$ws.Range("N7").Value = 30.370048
$ws.Range("O7").Value = 0.3206696780215441
$ws.Range("P7").Value = 0.3206696780215441
$ws.Range("Q7").Value = 242.0287254119538
$ws.Range("R7").Value = 2178.258528707584
$ws.Range("S7").Value = 0.07825655225147922
$ws.Range("T7").Value = 0.07825655225147925
$ws.Range("G8").Value = 47.93131266666666
$ws.Range("H8").Value = 143.793938
$ws.Range("I8").Value = 0.489259703040525
$ws.Range("J8").Value = 0.4892597030405251
$ws.Range("M8").Value = 0.8317113333333332
$ws.Range("N8").Value = 2.495134
$ws.Range("O8").Value = 0.0263454906755698
$ws.Range("P8").Value = 0.0263454906755698
$ws.Range("Q8").Value = 39.86501596641021
$ws.Range("R8").Value = 358.785143697692
$ws.Range("S8").Value = 0.0128897869443862
$ws.Range("T8").Value = 0.01288978694438621
$ws.Range("G9").Value = 47.93131266666666
$ws.Range("H9").Value = 143.793938
$ws.Range("I9").Value = 0.489259703040525
$ws.Range("J9").Value = 0.4892597030405251
$ws.Range("O9").Value = 0.6529848313028861
$ws.Range("P9").Value = 0.6529848313028862
$ws.Range("Q9").Value = 988.0723440027648
$ws.Range("R9").Value = 8892.651096024883
$ws.Range("S9").Value = 0.3194791646532174
$ws.Range("T9").Value = 0.3194791646532175
$ws.Range("G10").Value = 47.93131266666666
$ws.Range("H10").Value = 143.793938
$ws.Range("I10").Value = 0.489259703040525
$ws.Range("J10").Value = 0.4892597030405251
$ws.Range("M10").Value = 10.12334933333333
$ws.Range("N10").Value = 30.370048
$ws.Range("O10").Value = 0.3206696780215441
$ws.Range("P10").Value = 0.3206696780215441
$ws.Range("Q10").Value = 485.2254221298915
$ws.Range("R10").Value = 4367.028799169024
$ws.Range("S10").Value = 0.1568907514429214
$ws.Range("T10").Value = 0.1568907514429215
